$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 274.5
$ws.Range("J2").Value = 449
$ws.Range("L2").Value = 449
$ws.Range("N2").Value = -675

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 553.2222
$ws.Range("I58").Value = 432.25
$ws.Range("K58").Value = 1296.75
$ws.Range("M58").Value = -1146.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1344.2
$ws.Range("I112").Value = 1550
$ws.Range("J112").Value = 1326.3043
$ws.Range("K112").Value = 4650
$ws.Range("L112").Value = 3978.9129
$ws.Range("M112").Value = -3542
$ws.Range("N112").Value = -6194.9129

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 7180.04
$ws.Range("J137").Value = 5246.8667
$ws.Range("L137").Value = 15740.6001
$ws.Range("N137").Value = -20840.6001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7247.4
$ws.Range("I141").Value = 6362.385
$ws.Range("J141").Value = 13000
$ws.Range("K141").Value = 19087.155
$ws.Range("L141").Value = 39000
$ws.Range("M141").Value = -13907.155
$ws.Range("N141").Value = -49360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 52883404
$ws.Range("I5").Value = 6485287
$ws.Range("K5").Value = 6485287
$ws.Range("M5").Value = -6485175

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H30").Value = 3197.5
$ws.Range("I30").Value = 945
$ws.Range("K30").Value = 945
$ws.Range("M30").Value = -795

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13383.833
$ws.Range("I61").Value = 2357
$ws.Range("J61").Value = 35437.5
$ws.Range("K61").Value = 2357
$ws.Range("L61").Value = 35437.5
$ws.Range("M61").Value = -2145
$ws.Range("N61").Value = -35861.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 27786314
$ws.Range("I102").Value = 5435.4443
$ws.Range("J102").Value = 111128940
$ws.Range("K102").Value = 5435.4443
$ws.Range("L102").Value = 111128940
$ws.Range("M102").Value = -3813.4443
$ws.Range("N102").Value = -111132184

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5022.8335
$ws.Range("I122").Value = 4798.55
$ws.Range("J122").Value = 5471.4
$ws.Range("K122").Value = 14395.65
$ws.Range("L122").Value = 16414.2
$ws.Range("M122").Value = -11945.65
$ws.Range("N122").Value = -21314.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4573.3335
$ws.Range("I132").Value = 2888
$ws.Range("K132").Value = 8664
$ws.Range("M132").Value = -6134

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 13383.833
$ws.Range("I136").Value = 2357
$ws.Range("J136").Value = 35437.5
$ws.Range("K136").Value = 7071
$ws.Range("L136").Value = 106312.5
$ws.Range("M136").Value = -4521
$ws.Range("N136").Value = -111412.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 52883404
$ws.Range("I4").Value = 6485287
$ws.Range("K4").Value = 6485287
$ws.Range("M4").Value = -6485172

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5857.8823
$ws.Range("I20").Value = 4674
$ws.Range("J20").Value = 6910.222
$ws.Range("K20").Value = 4674
$ws.Range("L20").Value = 6910.222
$ws.Range("M20").Value = -4427
$ws.Range("N20").Value = -7404.222

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H36").Value = 1857
$ws.Range("I36").Value = 1857
$ws.Range("K36").Value = 1857
$ws.Range("M36").Value = -1323

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2418173
$ws.Range("I94").Value = 1993.6842
$ws.Range("K94").Value = 1993.6842
$ws.Range("M94").Value = -1542.6842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 18690
$ws.Range("I96").Value = 5305.1665
$ws.Range("K96").Value = 5305.1665
$ws.Range("M96").Value = -2559.1665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 40000
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 40000
$ws.Range("J137").Value = 40000
$ws.Range("L137").Value = 40000
$ws.Range("N137").Value = -50200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 58457.234
$ws.Range("J138").Value = 58457.234
$ws.Range("L138").Value = 58457.234
$ws.Range("N138").Value = -68737.234

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 40000
$ws.Range("J140").Value = 40000
$ws.Range("L140").Value = 40000
$ws.Range("N140").Value = -50360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 50000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 50000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 50000
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -50224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 255
$ws.Range("I7").Value = 129.75
$ws.Range("K7").Value = 129.75
$ws.Range("M7").Value = -16.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 44582.5
$ws.Range("J57").Value = 44582.5
$ws.Range("L57").Value = 44582.5
$ws.Range("N57").Value = -45702.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 6276
$ws.Range("I134").Value = 6445.2
$ws.Range("K134").Value = 19335.6
$ws.Range("M134").Value = -16800.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1745
$ws.Range("I14").Value = 1745
$ws.Range("K14").Value = 5235
$ws.Range("M14").Value = -5062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 732.5
$ws.Range("J109").Value = 375
$ws.Range("L109").Value = 1125
$ws.Range("N109").Value = -3205

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 613.4286
$ws.Range("I117").Value = 405
$ws.Range("J117").Value = 1134.5
$ws.Range("K117").Value = 1215
$ws.Range("L117").Value = 3403.5
$ws.Range("M117").Value = 2227
$ws.Range("N117").Value = -10287.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7569.25
$ws.Range("I70").Value = 7499
$ws.Range("K70").Value = 7499
$ws.Range("M70").Value = -7229

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7569.25
$ws.Range("I73").Value = 7499
$ws.Range("K73").Value = 7499
$ws.Range("M73").Value = -6563

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4915.75
$ws.Range("I102").Value = 4712.6665
$ws.Range("J102").Value = 5525
$ws.Range("K102").Value = 4712.6665
$ws.Range("L102").Value = 5525
$ws.Range("M102").Value = -3090.6665
$ws.Range("N102").Value = -8769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6939.6665
$ws.Range("I40").Value = 7666.1665
$ws.Range("J40").Value = 5486.6665
$ws.Range("K40").Value = 7666.1665
$ws.Range("L40").Value = 5486.6665
$ws.Range("M40").Value = -7530.1665
$ws.Range("N40").Value = -5758.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3083.7856
$ws.Range("I46").Value = 1166.6666
$ws.Range("K46").Value = 1166.6666
$ws.Range("M46").Value = -978.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1741.1875
$ws.Range("I61").Value = 1739.9286
$ws.Range("J61").Value = 1750
$ws.Range("K61").Value = 1739.9286
$ws.Range("L61").Value = 1750
$ws.Range("M61").Value = -1537.9286
$ws.Range("N61").Value = -2154

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 39995
$ws.Range("J63").Value = 39995
$ws.Range("L63").Value = 39995
$ws.Range("N63").Value = -41493

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H66").Value = 39995
$ws.Range("J66").Value = 39995
$ws.Range("L66").Value = 119985
$ws.Range("N66").Value = -127473

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1741.1875
$ws.Range("I113").Value = 1739.9286
$ws.Range("J113").Value = 1750
$ws.Range("K113").Value = 1739.9286
$ws.Range("L113").Value = 1750
$ws.Range("M113").Value = 430.0714
$ws.Range("N113").Value = -6090

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4121.875
$ws.Range("I136").Value = 3662.5
$ws.Range("J136").Value = 5500
$ws.Range("K136").Value = 10987.5
$ws.Range("L136").Value = 16500
$ws.Range("M136").Value = -8437.5
$ws.Range("N136").Value = -21600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 36666.5
$ws.Range("J70").Value = 36666.5
$ws.Range("L70").Value = 36666.5
$ws.Range("N70").Value = -37296.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 36666.5
$ws.Range("J73").Value = 36666.5
$ws.Range("L73").Value = 36666.5
$ws.Range("N73").Value = -38850.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1435.1538
$ws.Range("I122").Value = 1456.2222
$ws.Range("K122").Value = 4368.6666
$ws.Range("M122").Value = -1918.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2658.1667
$ws.Range("I126").Value = 2455.5557
$ws.Range("J126").Value = 3266
$ws.Range("K126").Value = 7366.6671
$ws.Range("L126").Value = 9798
$ws.Range("M126").Value = -4896.6671
$ws.Range("N126").Value = -14738

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4180.683
$ws.Range("I132").Value = 2718.9375
$ws.Range("K132").Value = 8156.8125
$ws.Range("M132").Value = -5626.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4000.0588
$ws.Range("I136").Value = 3937.875
$ws.Range("K136").Value = 11813.625
$ws.Range("M136").Value = -9263.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 52500
$ws.Range("J138").Value = 52500
$ws.Range("L138").Value = 52500
$ws.Range("N138").Value = -62780
